$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new (still blank) rows first: one at position 10 (pushes
#    the "smartphoneName" block and everything below it down by one row) and
#    one at position 17 (right after the "G6/g6" row, before the
#    "smartphoneRange" block, taking the earlier shift into account).
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(17).Insert()

$ws.Range("A10:C10").HorizontalAlignment = -4108
$ws.Range("A10:C10").VerticalAlignment = -4108
$ws.Range("C10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 30

$ws.Range("A17:C17").HorizontalAlignment = -4108
$ws.Range("A17:C17").VerticalAlignment = -4108
$ws.Range("C17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Fill in the new "6 / plus 6" entity row first, then the new "One Plus"
#    entity row, matching the order in which the entities were authored.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = 6
$ws.Range("C17").Value = "6`nplus 6"

$ws.Range("B10").Value = "One Plus"
$ws.Range("C10").Value = "one plus`noneplus"

# ---------------------------------------------------------------------------
# 3. Re-establish the merged "label" cells in column A so that they cover
#    the newly inserted rows as well.
# ---------------------------------------------------------------------------
$ws.Range("A3:A10").Merge()
$ws.Range("A11:A17").Merge()

# ---------------------------------------------------------------------------
# 4. Update the sheet selection to match the new state.
# ---------------------------------------------------------------------------
$ws.Range("C15:C16").Select()
